# Apply updated cryptocurrency price/volume data (and a few row content swaps)
# as described by the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.449.93"
$ws.Range("E2").Value = "  +12.38%  "
$ws.Range("D3").Value = "1.827.47"
$ws.Range("E3").Value = "  +9.06%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'230.18"
$ws.Range("E5").Value = "  +4.78%  "
$ws.Range("E6").Value = "  +7.99%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'31.54"
$ws.Range("E8").Value = "  +7.02%  "
$ws.Range("D9").Value = "'46.90"
$ws.Range("E9").Value = "  +5.79%  "
$ws.Range("D10").Value = "'0.288"
$ws.Range("E10").Value = "  +8.76%  "
$ws.Range("D11").Value = "'0.0679"
$ws.Range("E11").Value = "  +6.01%  "
$ws.Range("D12").Value = "'0.0932"
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("D13").Value = "2.091.34"
$ws.Range("E13").Value = "  +9.04%  "
$ws.Range("D14").Value = "1.841.23"
$ws.Range("E14").Value = "  +9.85%  "
$ws.Range("D15").Value = "'0.656"
$ws.Range("E15").Value = "  +8.17%  "
$ws.Range("D16").Value = "34.396.59"
$ws.Range("E16").Value = "  +12.11%  "
$ws.Range("D17").Value = "'10.30"
$ws.Range("E17").Value = "  +3.85%  "
$ws.Range("D18").Value = "'4.31"
$ws.Range("E18").Value = "  +7.20%  "
$ws.Range("D19").Value = "'70.49"
$ws.Range("E19").Value = "  +6.37%  "
$ws.Range("D20").Value = "'258.23"
$ws.Range("E20").Value = "  +6.53%  "
$ws.Range("D21").Value = "0.0₃0759"
$ws.Range("E21").Value = "  +5.19%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "'10.64"
$ws.Range("E23").Value = "  +6.56%  "
$ws.Range("D24").Value = "'4.34"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("D25").Value = "'2.23"
$ws.Range("E25").Value = "  +3.68%  "
$ws.Range("D26").Value = "'159.56"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "'16.81"
$ws.Range("E27").Value = "  +6.33%  "
$ws.Range("E28").Value = "  +4.46%  "
$ws.Range("D29").Value = "'7.18"
$ws.Range("E29").Value = "  +7.43%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "'3.89"
$ws.Range("E31").Value = "  +12.21%  "
$ws.Range("D32").Value = "'0.0524"
$ws.Range("E32").Value = "  +6.35%  "
$ws.Range("D33").Value = "'1.21"
$ws.Range("E33").Value = "  +6.07%  "
$ws.Range("D34").Value = "'3.60"
$ws.Range("E34").Value = "  +8.23%  "
$ws.Range("D35").Value = "1.548.43"
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("D36").Value = "'1.80"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").Value = "'1.07"
$ws.Range("E37").Value = "  +4.89%  "
$ws.Range("D38").Value = "'0.639"
$ws.Range("E38").Value = "  +6.86%  "
$ws.Range("D39").Value = "'0.0191"
$ws.Range("E39").Value = "  +7.27%  "
$ws.Range("D40").Value = "'84.73"
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("D41").Value = "'2.82"
$ws.Range("E41").Value = "  +5.15%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "'2.36"
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.915"
$ws.Range("E43").Value = "  +9.09%  "
$ws.Range("D44").Value = "'2.14"
$ws.Range("E44").Value = "  +6.04%  "
$ws.Range("E45").Value = "  +5.78%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.08"
$ws.Range("E46").Value = "  +6.17%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.976.28"
$ws.Range("E47").Value = "  +9.01%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'5.84"
$ws.Range("E48").Value = "  +5.00%  "
$ws.Range("D49").Value = "'12.27"
$ws.Range("E49").Value = "  +19.07%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'51.84"
$ws.Range("E51").Value = "  +2.83%  "
